$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: date changed from 2013-09-10 (41527) to 2013-10-10 (41557)
$ws.Range("A6").Value = 41557

# Row 7: new date entry (2013-10-11) with matching date style copied from A6,
# plus the hours worked that day (3:20 -> 0.1388888888888889)
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)   # xlPasteFormats, so A7 gets style s="1" like the other date cells
$ws.Range("A7").Value = 41558
$ws.Range("B7").Value = 0.1388888888888889

# Reflect the new rows in the sheet's selection (B4:B7 instead of just B7)
$ws.Range("B4:B7").Select()
